# Update "想去人数" (interested-count) figures to the latest scrape snapshot.
# Sheet 1 = 展览 (Exhibitions), Sheet 2 = 演出 (Performances),
# Sheet 3 = 本地生活 (Local Life, unchanged), Sheet 4 = 全部类型 (All types, aggregate view).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsShow    = $wb.Worksheets.Item(2)   # 演出
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# 展览 (sheet1)
$wsExhibit.Range("F5").Value  = 1081
$wsExhibit.Range("F6").Value  = 171
$wsExhibit.Range("F8").Value  = 212
$wsExhibit.Range("F9").Value  = 387
$wsExhibit.Range("F11").Value = 10
$wsExhibit.Range("F14").Value = 153
$wsExhibit.Range("F15").Value = 12532
$wsExhibit.Range("F16").Value = 145
$wsExhibit.Range("F17").Value = 5502

# 演出 (sheet2)
$wsShow.Range("F2").Value = 123
$wsShow.Range("F4").Value = 2

# 全部类型 (sheet4) - same underlying rows aggregated, different row numbers
$wsAll.Range("F4").Value  = 123
$wsAll.Range("F7").Value  = 1081
$wsAll.Range("F8").Value  = 171
$wsAll.Range("F10").Value = 212
$wsAll.Range("F11").Value = 387
$wsAll.Range("F13").Value = 10
$wsAll.Range("F16").Value = 153
$wsAll.Range("F17").Value = 12532
$wsAll.Range("F18").Value = 2
$wsAll.Range("F20").Value = 145
$wsAll.Range("F21").Value = 5502
